$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (serial date, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
# covering the update "fino a 6 gennaio 2022"
$newRows = @(
    @(465, 44539, 5, 9, 394.045534150613),
    @(466, 44540, 1, 9, 394.045534150613),
    @(467, 44541, 0, 9, 394.045534150613),
    @(468, 44542, 1, 8, 350.2626970227671),
    @(469, 44543, 1, 9, 394.045534150613),
    @(470, 44544, 2, 11, 481.6112084063047),
    @(471, 44545, 0, 10, 437.8283712784589),
    @(472, 44546, 0, 5, 218.9141856392294),
    @(473, 44547, 2, 6, 262.6970227670753),
    @(474, 44548, 1, 7, 306.4798598949212),
    @(475, 44550, 2, 8, 350.2626970227671),
    @(476, 44551, 0, 7, 306.4798598949212),
    @(477, 44552, 0, 5, 218.9141856392294),
    @(478, 44553, 0, 5, 218.9141856392294),
    @(479, 44554, 0, 5, 218.9141856392294),
    @(480, 44555, 0, 3, 131.3485113835376),
    @(481, 44556, 0, 2, 87.56567425569177),
    @(482, 44557, 0, 0, 0),
    @(483, 44558, 0, 0, 0),
    @(484, 44559, 1, 1, 43.78283712784589),
    @(485, 44560, 1, 2, 87.56567425569177),
    @(486, 44561, 1, 3, 131.3485113835376),
    @(487, 44562, 4, 7, 306.4798598949212),
    @(488, 44563, 0, 7, 306.4798598949212),
    @(489, 44564, 1, 8, 350.2626970227671),
    @(490, 44565, 1, 9, 394.045534150613),
    @(491, 44566, 0, 8, 350.2626970227671)
)

$firstNewRow = 465
$lastNewRow = 491

# Copy the formatting (date number format / alignment / border) of the last
# existing row down onto the newly appended rows before filling in values.
$ws.Range("A464").Copy()
$ws.Range("A" + $firstNewRow + ":A" + $lastNewRow).PasteSpecial(-4122)

foreach ($item in $newRows) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
    $ws.Cells.Item($r, 4).Value = $item[4]
}

$excel.CutCopyMode = 0
